# Weekly driver report update for 2025-04-20
# Refresh / re-sort the "Good Drivers" table (rows 12-17) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-VintageDate($cellRef, $dateText) {
    # Writing an ISO-ish date string through .Value auto-converts it to a
    # real date serial + date number-format. The source report stores the
    # "Driver Vintage" column as plain text, right-aligned like the numeric
    # columns next to it (style index 4: General format, right alignment).
    # Force text storage, write the literal text, then restore the shared
    # "General / right aligned" look so the cell lands on the same style
    # as its neighbours instead of minting a new date-formatted style.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $dateText
    $cell.ClearFormats()
    $cell.HorizontalAlignment = -4152   # xlRight
}

function Clear-VintageDate($cellRef) {
    $cell = $ws.Range($cellRef)
    $cell.Value = ""
}

# Row 12: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
Set-VintageDate "E12" "2024-11-10"

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
Set-VintageDate "E13" "2021-08-18"

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
Set-VintageDate "E14" "2021-04-27"

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
Set-VintageDate "E15" "2020-08-05"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("D16").Value = 100
Set-VintageDate "E16" "2020-01-06"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
$ws.Range("D17").Value = 100
Clear-VintageDate "E17"
